# "Finishing the Location check"
# - Sheet1: correct two people's Location from Szeged -> London (B5, B7),
#   fix D5's project value, add "SHOT 2" to the June-column autofilter
#   (which reveals rows 2, 3 and 11), make Sheet1 the active/selected tab
#   with selection C17, and set the page orientation.
# - projectDetails: fix the "SHOT" row's allowed Location (C3 -> Szeged)
#   and its max-location-allowed count (D5 -> 2); selection moves to D3
#   and the tab is no longer the one left "active" in the saved view.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("projectDetails")

# Row 7 is (and stays) hidden; temporarily reveal it so the value write
# below doesn't disturb its row height, then restore it before the
# AutoFilter call recomputes every row's visibility.
$ws1.Rows.Item(7).Hidden = $false

$ws1.Range("B5").Value = "London"
$ws1.Range("D5").Value = "ROOMS AND RATES"
$ws1.Range("B7").Value = "London"

$ws1.Rows.Item(7).Hidden = $true

# Add "SHOT 2" alongside the existing "SHOT" criterion on the June column
# (field 5 = column E). This also re-evaluates hidden/visible rows for
# every record against the new multi-value filter.
$null = $ws1.Range("A1:E14").AutoFilter(5, @("SHOT", "SHOT 2"), 7)

$ws1.PageSetup.Orientation = 1

# projectDetails data fixes
$ws2.Range("C3").Value = "Szeged"
$ws2.Range("D5").Value = 2

# Move the workbook window geometry to the saved-default state.
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 0
$win.Width = 28800
$win.Height = 18000

# Selections: Sheet1 becomes the active tab with C17 selected; the
# projectDetails sheet keeps D3 selected but is no longer the active tab.
$null = $ws2.Range("D3").Select()
$null = $ws1.Range("C17").Select()
